$d = $word.ActiveDocument

# Update the order date: 15/08/2014 -> 29/09/2014
$d.Content.Find.Execute("15/08/2014", $true, $false, $false, $false, $false,
                         $true, 1, $false, "29/09/2014", 2)

# Update the board size: 138.4 x 86.2 mm -> 138.55 x 86.23 mm
$d.Content.Find.Execute("138.4 x 86.2 mm", $true, $false, $false, $false, $false,
                         $true, 1, $false, "138.55 x 86.23 mm", 2)
